$wb = $excel.ActiveWorkbook

# --- OFF sheet (row 2, "H") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 465
$wsOff.Range("C2").Value = 295
$wsOff.Range("D2").Value = 131
$wsOff.Range("E2").Value = 55
$wsOff.Range("F2").Value = 22

# --- DEF sheet (row 2, "H") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 469
$wsDef.Range("C2").Value = 317
$wsDef.Range("D2").Value = 133
$wsDef.Range("E2").Value = 57
$wsDef.Range("G2").Value = 11
